# Atualização de bases das ligas, do dia: 02-05-2024 às 20:28
#
# Rows 22-24 rotate cyclically (row22 <- row23, row23 <- row24, row24 <- row22)
# and rows 33-34 swap with each other, keeping column A (the running index)
# fixed on each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AB hold the data that moves; column A stays put.
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB")

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$row").Value2
    }
    return $vals
}

function Set-RowValues($ws, $row, $cols, $vals) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value2 = $vals[$c]
    }
}

# Snapshot current (pre-edit) values for the affected rows.
$row22 = Get-RowValues $ws 22 $cols
$row23 = Get-RowValues $ws 23 $cols
$row24 = Get-RowValues $ws 24 $cols
$row33 = Get-RowValues $ws 33 $cols
$row34 = Get-RowValues $ws 34 $cols

# Cyclic rotation: 22 <- 23, 23 <- 24, 24 <- 22
Set-RowValues $ws 22 $cols $row23
Set-RowValues $ws 23 $cols $row24
Set-RowValues $ws 24 $cols $row22

# Swap rows 33 and 34
Set-RowValues $ws 33 $cols $row34
Set-RowValues $ws 34 $cols $row33
